# Add two new slides ("Results" and "Conclusion") at the end of the
# deck, using the "Title and Content" layout (the same layout already
# used by slides 2-4: ppt/slideLayouts/slideLayout2.xml -> layout index 2).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 5 - "Results"
# ---------------------------------------------------------------
$s5 = $p.Slides.Add(5, 2)

$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Results"

$body5 = $s5.Shapes.Item(2).TextFrame.TextRange
$body5.Text = "Providence has more venues than Hartford (512 vs 460)"
$null = $body5.InsertAfter("`rFederal Hill has the most venues in Providence")
$null = $body5.InsertAfter("`rDowntown has the most venues in Hartford")
$null = $body5.InsertAfter("`rMost of venues in Federal Hill are restaurants, bars")
$null = $body5.InsertAfter("`rVenues in Downtown Hartford are more diverse includes: banks, concert hall, gyms…")
$null = $body5.InsertAfter("`rBars, Restaurants, Bakeries, and Diners are very common in Providence")
$null = $body5.InsertAfter("`rHartford is more about commercial")

# ---------------------------------------------------------------
# Slide 6 - "Conclusion"
# ---------------------------------------------------------------
$s6 = $p.Slides.Add(6, 2)

$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Conclusion"

$body6 = $s6.Shapes.Item(2).TextFrame.TextRange
$body6.Text = "For travelers:"
$null = $body6.InsertAfter("`rProvidence may more suitable to spend a good amount time of enjoy the local food")
$null = $body6.InsertAfter("`rHartford has less food but more commercial")
$null = $body6.InsertAfter("`r")
$null = $body6.InsertAfter("`rFor business owners:")
$null = $body6.InsertAfter("`rProvidence is good place to open small business related to food and travel industry")
$null = $body6.InsertAfter("`rHartford is better place to involve in insurance industry")
$null = $body6.InsertAfter("`r")

# The trailing blank paragraph has its bullet cleared (8 paragraphs
# total were built above; the COM "Paragraphs().Count" undercounts a
# final empty paragraph by one, so address it by its known, fixed
# position instead of relying on that count).
$lastPara = $s6.Shapes.Item(2).TextFrame.TextRange.Paragraphs(8, 1)
$lastPara.ParagraphFormat.Bullet.Visible = 0
